# Update "想去人数" (people interested) counts in column F across sheets.
# Mirrors the commit "Update gh-pages to output generated at 456a3b4":
#   展览 (Exhibitions), 演出 (Performances), 本地生活 (Local Life) each get
#   fresh counts, and 全部类型 (All Types) is the aggregate sheet that
#   repeats the same rows, so it gets the matching updates too.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 389
$ws.Range("F5").Value = 1326
$ws.Range("F6").Value = 232
$ws.Range("F7").Value = 2520
$ws.Range("F8").Value = 925
$ws.Range("F9").Value = 18713
$ws.Range("F10").Value = 56
$ws.Range("F11").Value = 1943
$ws.Range("F12").Value = 672
$ws.Range("F15").Value = 609
$ws.Range("F18").Value = 74
$ws.Range("F19").Value = 322
$ws.Range("F20").Value = 35
$ws.Range("F21").Value = 4
$ws.Range("F23").Value = 113

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 111

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5899
$ws.Range("F3").Value = 576
$ws.Range("F4").Value = 560

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5899
$ws.Range("F4").Value = 576
$ws.Range("F5").Value = 560
$ws.Range("F6").Value = 389
$ws.Range("F10").Value = 1326
$ws.Range("F12").Value = 232
$ws.Range("F15").Value = 2520
$ws.Range("F16").Value = 925
$ws.Range("F17").Value = 18713
$ws.Range("F20").Value = 56
$ws.Range("F24").Value = 1943
$ws.Range("F25").Value = 672
$ws.Range("F26").Value = 111
$ws.Range("F28").Value = 609
$ws.Range("F32").Value = 74
$ws.Range("F35").Value = 322
$ws.Range("F36").Value = 35
$ws.Range("F38").Value = 4
$ws.Range("F41").Value = 113
